$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LogIn")

# Update the email value in A2 (shared string changes from navin123@gmail.com to parag123@gmail.com)
$ws.Range("A2").Value = "parag123@gmail.com"

# Move the active selection from C8 to B6
$ws.Activate()
$ws.Range("B6").Select()
